# Updated cryptos list on Mon Nov  4 23:52:47 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.737.90"
$ws.Range("E2").Value = "  -1.66%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.395.66"
$ws.Range("E3").Value = "  -2.55%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'551.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.17%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'157.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.39%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.06%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.503"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +3.97%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  -1.39%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -2.64%  "

# Row 12 - Toncoin
$ws.Range("D12").Value = "'4.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.58%  "

# Row 13 - WrappedBTC
$ws.Range("D13").Value = "67.655.26"
$ws.Range("E13").Value = "  -1.50%  "

# Row 14 - ShibaInu
$ws.Range("D14").Value = "'0.0000168"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.63%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "'22.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.76%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'10.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.87%  "

# Row 17 - BitcoinCash
$ws.Range("D17").Value = "'329.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.50%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "'6.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.05%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -1.25%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  -0.05%  "

# Row 21 - SuiNetwork
$ws.Range("D21").Value = "'1.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.24%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "'65.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.29%  "

# Row 23 - NEARProtocol
$ws.Range("D23").Value = "'3.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.25%  "

# Row 24 - Aptos
$ws.Range("E24").Value = "  -2.17%  "

# Row 25 - PEPE
$ws.Range("D25").Value = "0.0₃0791"
$ws.Range("E25").Value = "  -3.53%  "

# Row 26 - InternetComputer(DFINITY)
$ws.Range("D26").Value = "'7.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.36%  "

# Row 27 - FirstDigitalUSD
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

# Row 28 - Bittensor
$ws.Range("D28").Value = "'418.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.85%  "

# Row 29 - Fetch.AI
$ws.Range("E29").Value = "  -1.50%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -2.20%  "

# Row 31 - Monero
$ws.Range("D31").Value = "'157.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.40%  "

# Row 32 - WhiteBITCoin
$ws.Range("D32").Value = "'18.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.42%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  -0.04%  "

# Row 34 - EthereumClassic
$ws.Range("D34").Value = "'17.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.32%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  -4.41%  "

# Row 36 - PolygonEcosystemToken
$ws.Range("E36").Value = "  -3.27%  "

# Row 37 - RenderToken
$ws.Range("D37").Value = "'4.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.50%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  -2.01%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "'1.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.70%  "

# Row 42 - dogwifhat
$ws.Range("D42").Value = "'1.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.98%  "

# Row 43 - Cronos
$ws.Range("D43").Value = "'0.0704"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.17%  "

# Row 44 - ARBITRUM
$ws.Range("D44").Value = "'0.473"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.95%  "

# Row 45 - Mantle
$ws.Range("D45").Value = "'0.551"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.58%  "

# Row 46 - Stellar
$ws.Range("D46").Value = "'0.0911"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.50%  "

# Row 47 - BitgetToken
$ws.Range("D47").Value = "'1.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.12%  "

# Row 48 - Optimism
$ws.Range("D48").Value = "'1.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.47%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "'16.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.52%  "

# Row 50 - BabyDogeCoin
$ws.Range("D50").Value = "0.0⁦0201"
$ws.Range("E50").Value = "  -7.22%  "

# Row 51 - Hedera
$ws.Range("D51").Value = "'0.0425"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.89%  "

# Row 40 & 41 - Aave / Filecoin swap positions with new data
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'3.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.30%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'128.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.58%  "

